$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Thursday time in the header row (D1): was "14:15 - 16:00", now "12:15 - 14:00"
$ws.Range("D1").Value = "Thursday 12:15 - 14:00"

# Move the active selection to D2
$ws.Range("D2").Select()
